# Add a new "trades" worksheet as the last sheet in the workbook, mirroring the
# header style already used on the other sheets (e.g. "stocks"!A1).

$wb = $excel.ActiveWorkbook

# Remember the currently active sheet so we can restore it afterwards (adding
# a sheet normally activates it, but we don't want to change the workbook's
# selected/active tab as part of this edit).
$originalActiveSheet = $wb.Worksheets.Item(1)

# Grab the style used for header cells elsewhere in the workbook so the new
# sheet's headers look consistent (bold, bordered, centered) and reuse the
# existing style definition instead of creating a new one.
$stocksSheet = $wb.Worksheets.Item("stocks")
$headerStyleSource = $stocksSheet.Range("A1")

# Insert the new sheet after the current last sheet so it lands at the end
# (matching sheetId="4" / position 4 in the target workbook).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tradesSheet = $wb.Worksheets.Add($null, $lastSheet)
$tradesSheet.Name = "trades"

# Populate the header row starting at column B (matches target layout B1:F1).
$tradesSheet.Range("B1").Value = "date"
$tradesSheet.Range("C1").Value = "ticker"
$tradesSheet.Range("D1").Value = "buy_sell"
$tradesSheet.Range("E1").Value = "shares"
$tradesSheet.Range("F1").Value = "value"

# Copy the header formatting (font/border/alignment) onto the new header row.
$headerStyleSource.Copy()
$tradesSheet.Range("B1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the originally active sheet/tab.
$originalActiveSheet.Activate()
